$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new columns (P, Q) to the header row, copying the format
#     from the last existing header cell (O1) so the same style index
#     ("s=1") is reused instead of creating a brand-new style entry.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Fill in the new P and Q columns for the data rows (2-25) with 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 16).Value = 2   # P
    $ws.Cells.Item($r, 17).Value = 2   # Q
}

# --- Update the swapped values in columns I, K, M, O for rows 2-25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
}

Write-Output "edit applied"
